$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "5(1)"
$ws.Range("C12").Value = "5(2)"
$ws.Range("D9").Value = "5(3)"
$ws.Range("E14").Value = "5(4)"
$ws.Range("F11").Value = "5(5)"
$ws.Range("H13").Value = "5(6)"
$ws.Range("J15").Value = "5(7)"
$ws.Range("B11").Value = "6(1)"
$ws.Range("C8").Value = "6(2)"
$ws.Range("D13").Value = "6(3)"
$ws.Range("E10").Value = "6(4)"
$ws.Range("F15").Value = "6(5)"
$ws.Range("G12").Value = "6(6)"
$ws.Range("I14").Value = "6(7)"

$ws.Range("I15").Select()
